$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original (hyperlink) cell style so it can be restored after
# re-creating the hyperlinks below (Hyperlinks.Add() re-applies the builtin
# "Hyperlink" style as a brand new cellXf, which we don't want to change).
$linkStyle = $ws.Range("J2").Style

# --- Update row 3 ("Ambiente"/"URL" columns) to point at the new,
#     non "i-" prefixed preproduccion host -------------------------------
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# The runtime does not support editing a Hyperlink's Address in place, so
# clear every hyperlink on the sheet and recreate them, fixing only the
# one that changed (B3); the mailto: links on J2/J3 are restored as-is.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:aseguradosgw@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:aseguradosgw@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do") | Out-Null

# Restore original cell styles (Hyperlinks.Add() reassigns a fresh style).
$ws.Range("J2").Style = $linkStyle
$ws.Range("J3").Style = $linkStyle
$ws.Range("B3").Style = $linkStyle

# --- Update the saved selection / active cell ---------------------------
$ws.Range("B4").Select()
